$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '52.078.26'
$ws.Range('E2').Value = '  +1.32%  '
$ws.Range('D3').Value = '2.881.47'
$ws.Range('E3').Value = '  +3.74%  '
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').Value = "'352.18"
$ws.Range('E5').Value = '  -0.12%  '
$ws.Range('D6').Value = "'111.56"
$ws.Range('E6').Value = '  +3.32%  '
$ws.Range('D7').Value = "'0.559"
$ws.Range('E7').Value = '  +1.76%  '
$ws.Range('E8').Value = '  +0.05%  '
$ws.Range('D9').Value = "'0.621"
$ws.Range('E9').Value = '  +0.81%  '
$ws.Range('D10').Value = "'40.08"
$ws.Range('E10').Value = '  +2.46%  '
$ws.Range('D11').Value = "'0.0863"
$ws.Range('E11').Value = '  +3.61%  '
$ws.Range('D12').Value = "'0.136"
$ws.Range('E12').Value = '  +0.47%  '
$ws.Range('D13').Value = "'20.05"
$ws.Range('E13').Value = '  +0.99%  '
$ws.Range('D14').Value = "'7.82"
$ws.Range('E14').Value = '  +1.04%  '
$ws.Range('D15').Value = '3.337.27'
$ws.Range('E15').Value = '  +3.97%  '
$ws.Range('B16').Value = 'Polygon'
$ws.Range('C16').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D16').Value = "'0.995"
$ws.Range('E16').Value = '  +7.86%  '
$ws.Range('B17').Value = 'WrappedEther'
$ws.Range('C17').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D17').Value = '2.878.55'
$ws.Range('E17').Value = '  +3.01%  '
$ws.Range('D18').Value = '52.093.90'
$ws.Range('E18').Value = '  +1.44%  '
$ws.Range('E19').Value = '  +8.09%  '
$ws.Range('E20').Value = '  +0.17%  '
$ws.Range('D21').Value = "'13.88"
$ws.Range('E21').Value = '  +4.06%  '
$ws.Range('E22').Value = '  +1.80%  '
$ws.Range('D23').Value = "'70.98"
$ws.Range('E23').Value = '  +0.77%  '
$ws.Range('D24').Value = "'270.66"
$ws.Range('E24').Value = '  +1.75%  '
$ws.Range('E25').Value = '  +0.91%  '
$ws.Range('D26').Value = "'26.31"
$ws.Range('E26').Value = '  +1.97%  '
$ws.Range('D27').Value = "'0.998"
$ws.Range('E27').Value = '  -0.14%  '
$ws.Range('D28').Value = "'0.163"
$ws.Range('E28').Value = '  +0.18%  '
$ws.Range('D29').Value = "'10.54"
$ws.Range('E29').Value = '  +2.81%  '
$ws.Range('D30').Value = "'38.73"
$ws.Range('E30').Value = '  +4.59%  '
$ws.Range('E31').Value = '  +0.75%  '
$ws.Range('D32').Value = "'6.43"
$ws.Range('E32').Value = '  +3.29%  '
$ws.Range('D33').Value = "'53.34"
$ws.Range('E33').Value = '  +3.05%  '
$ws.Range('B34').Value = 'Hedera'
$ws.Range('C34').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D34').Value = "'0.0940"
$ws.Range('E34').Value = '  +12.15%  '
$ws.Range('B35').Value = 'RenderToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D35').Value = "'5.91"
$ws.Range('E35').Value = '  +4.88%  '
$ws.Range('E36').Value = '  +3.91%  '
$ws.Range('E37').Value = '  -0.06%  '
$ws.Range('E38').Value = '  +6.77%  '
$ws.Range('D39').Value = "'18.59"
$ws.Range('E39').Value = '  +0.82%  '
$ws.Range('E40').Value = '  +3.69%  '
$ws.Range('D41').Value = "'2.64"
$ws.Range('E41').Value = '  +6.46%  '
$ws.Range('E42').Value = '  +2.93%  '
$ws.Range('D43').Value = "'22.47"
$ws.Range('E43').Value = '  +3.09%  '
$ws.Range('D44').Value = "'121.69"
$ws.Range('E44').Value = '  +1.36%  '
$ws.Range('E45').Value = '  +1.39%  '
$ws.Range('D46').Value = "'3.58"
$ws.Range('E46').Value = '  +7.17%  '
$ws.Range('D47').Value = '2.200.96'
$ws.Range('E47').Value = '  +3.30%  '
$ws.Range('D48').Value = "'2.49"
$ws.Range('E48').Value = '  +6.57%  '
$ws.Range('D49').Value = "'0.271"
$ws.Range('E49').Value = '  +19.93%  '
$ws.Range('D50').Value = "'0.951"
$ws.Range('E50').Value = '  +6.92%  '
$ws.Range('D51').Value = "'5.51"
$ws.Range('E51').Value = '  +1.01%  '
